# updated sustainable tourism goal
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column F width: widen from "best fit" 15.33 to a fixed 21.33 chars ---
$ws.Columns.Item(6).ColumnWidth = 20.5

# --- Row 2: revise visitor-spending figures upward ---
$ws.Range("E2").NumberFormat = "General"
$ws.Range("E2").Value = 160.4
$ws.Range("F2").NumberFormat = "General"
$ws.Range("F2").Value = 1839

# --- Row 3: align F3's number format with the other "General" cells ---
$ws.Range("F3").NumberFormat = "General"

# --- Add the 2010 data row ---
$ws.Range("B7").Value = 2010
$ws.Range("C7").Value = 1378921
$ws.Range("D7").Value = 7284769
$ws.Range("E7").Value = 145.1
$ws.Range("F7").Value = 1345.7

# --- Move the active selection ---
$ws.Range("G14").Select()
